# The deck currently has the "Integral" (Red Violet) design applied to its
# single slide master/theme. The author switched the presentation's design
# back to the built-in default "Office Theme" colour palette via the
# PowerPoint Design gallery. Reproduce that by rewriting the 12 theme colour
# slots (clrScheme) on the presentation's theme through the Slide Master's
# ThemeColorScheme, which is the supported COM surface for editing theme
# colours.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$colors = $master.Theme.ThemeColorScheme

# Office Theme colour scheme (RGB() uses 0xBBGGRR ordering)
$colors.Item(1).RGB  = 0x000000   # dk1
$colors.Item(2).RGB  = 0xFFFFFF   # lt1
$colors.Item(3).RGB  = 0x6A5444   # dk2      (44546A)
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      (E7E6E6)
$colors.Item(5).RGB  = 0xD59B5B   # accent1  (5B9BD5)
$colors.Item(6).RGB  = 0x317DED   # accent2  (ED7D31)
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  (A5A5A5)
$colors.Item(8).RGB  = 0x00C0FF   # accent4  (FFC000)
$colors.Item(9).RGB  = 0xC47244   # accent5  (4472C4)
$colors.Item(10).RGB = 0x47AD70   # accent6  (70AD47)
$colors.Item(11).RGB = 0xC16305   # hlink    (0563C1)
$colors.Item(12).RGB = 0x724F95   # folHlink (954F72)
